$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2 to whole numbers
$ws.Range("Q2").Value = 700829
$ws.Range("R2").Value = 7327307

# Remove the time cells Z2 (Starttid) and AB2 (Sluttid)
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
